$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number, week date range) ---
$ws.Range("A8").Value = "Volume 31   Number  17"
$ws.Range("C9").Value = "Report Covering the Week  4/22/2024  Through  4/28/2024"

# --- Crime-data table updates (rows 15-31) ---
# Row 15
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("E15").Value = -100
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G15").Value = 1
$ws.Range("G15").NumberFormat = "#,##0"
$ws.Range("H15").Value = -100
$ws.Range("H15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J15").Value = 1
$ws.Range("J15").NumberFormat = "#,##0"
$ws.Range("K15").Value = 300
$ws.Range("K15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L15").Value = 33.333333333333

# Row 16
$ws.Range("C16").Value = 6
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "***.*"
$ws.Range("E16").NumberFormat = "General"
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 4
$ws.Range("H16").Value = 200
$ws.Range("I16").Value = 27
$ws.Range("K16").Value = -3.571428571428
$ws.Range("L16").Value = 28.571428571428
$ws.Range("M16").Value = 28.571428571428
$ws.Range("N16").Value = -85.483870967741

# Row 17
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 13
$ws.Range("G17").Value = 8
$ws.Range("H17").Value = 62.5
$ws.Range("I17").Value = 42
$ws.Range("J17").Value = 37
$ws.Range("K17").Value = 13.513513513513
$ws.Range("L17").Value = 5
$ws.Range("N17").Value = -20.754716981132

# Row 18
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = 10
$ws.Range("I18").Value = 35
$ws.Range("J18").Value = 45
$ws.Range("K18").Value = -22.222222222222
$ws.Range("L18").Value = -46.153846153846
$ws.Range("M18").Value = -12.5
$ws.Range("N18").Value = -90.716180371352

# Row 19
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 20
$ws.Range("E19").Value = -50
$ws.Range("F19").Value = 37
$ws.Range("G19").Value = 81
$ws.Range("H19").Value = -54.320987654321
$ws.Range("I19").Value = 189
$ws.Range("J19").Value = 236
$ws.Range("K19").Value = -19.915254237288
$ws.Range("L19").Value = -3.076923076923
$ws.Range("M19").Value = -14.479638009049
$ws.Range("N19").Value = -69.614147909967

# Row 20
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "0"
$ws.Range("C20").NumberFormat = "General"
$ws.Range("E20").Value = -100
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = -66.666666666666
$ws.Range("J20").Value = 16
$ws.Range("K20").Value = -50
$ws.Range("L20").Value = -11.111111111111
$ws.Range("M20").Value = 60
$ws.Range("N20").Value = -96.226415094339

# Row 21
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = -25.925925925925
$ws.Range("F21").Value = 74
$ws.Range("G21").Value = 107
$ws.Range("H21").Value = -30.841121495327
$ws.Range("I21").Value = 306
$ws.Range("J21").Value = 363
$ws.Range("K21").Value = -15.702479338843
$ws.Range("L21").Value = -8.108108108108
$ws.Range("M21").Value = -1.6077170418
$ws.Range("N21").Value = -78.954607977991

# Row 22
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("C22").NumberFormat = "General"

# Row 24
$ws.Range("C24").Value = 15
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = -21.052631578947
$ws.Range("F24").Value = 66
$ws.Range("G24").Value = 83
$ws.Range("H24").Value = -20.481927710843
$ws.Range("I24").Value = 313
$ws.Range("J24").Value = 303
$ws.Range("K24").Value = 3.300330033003
$ws.Range("L24").Value = -16.085790884718
$ws.Range("M24").Value = 51.941747572815

# Row 25
$ws.Range("C25").Value = 17
$ws.Range("D25").Value = 16
$ws.Range("E25").Value = 6.25
$ws.Range("F25").Value = 60
$ws.Range("G25").Value = 68
$ws.Range("H25").Value = -11.764705882352
$ws.Range("I25").Value = 264
$ws.Range("J25").Value = 244
$ws.Range("K25").Value = 8.196721311475
$ws.Range("L25").Value = -15.112540192926

# Row 26
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = -60
$ws.Range("F26").Value = 14
$ws.Range("G26").Value = 19
$ws.Range("H26").Value = -26.315789473684
$ws.Range("I26").Value = 69
$ws.Range("J26").Value = 79
$ws.Range("K26").Value = -12.658227848101
$ws.Range("L26").Value = -9.210526315789
$ws.Range("M26").Value = -2.81690140845

# Row 27
$ws.Range("D27").Value = 1
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("E27").Value = -100
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G27").Value = 1
$ws.Range("G27").NumberFormat = "#,##0"
$ws.Range("H27").Value = -100
$ws.Range("H27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J27").Value = 3
$ws.Range("K27").Value = 66.666666666666
$ws.Range("L27").Value = 66.666666666666

# Row 28
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 6
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 19
$ws.Range("J28").Value = 19
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 5.555555555555

# Row 31
$ws.Range("D31").Value = 1
$ws.Range("D31").NumberFormat = "#,##0"
$ws.Range("E31").Value = -100
$ws.Range("E31").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J31").Value = 3
$ws.Range("K31").Value = 33.333333333333
